$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

for ($i = 1; $i -le 3; $i++) {
    $ftr = $sec.Footers.Item($i)
    for ($j = $ftr.Shapes.Count; $j -ge 1; $j--) {
        $shp = $ftr.Shapes.Item($j)
        $shp.Delete()
    }
}
